$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.165.34'
$ws.Range('E2').Value = '  -0.10%  '
$ws.Range('D3').Value = '1.842.36'
$ws.Range('E3').Value = '  -0.31%  '
$ws.Range('D5').Value = "'241.63"
$ws.Range('E5').Value = '  -1.61%  '
$ws.Range('D6').Value = "'0.6879"
$ws.Range('E6').Value = '  -1.89%  '
$ws.Range('D7').Value = "'0.9995"
$ws.Range('E7').Value = '  -0.10%  '
$ws.Range('D8').Value = "'0.3016"
$ws.Range('E8').Value = '  -1.65%  '
$ws.Range('D9').Value = "'0.07464"
$ws.Range('E9').Value = '  -3.29%  '
$ws.Range('D10').Value = "'23.14"
$ws.Range('E10').Value = '  -1.94%  '
$ws.Range('D11').Value = "'0.07657"
$ws.Range('E11').Value = '  -1.93%  '
$ws.Range('D12').Value = '1.840.25'
$ws.Range('E12').Value = '  -0.63%  '
$ws.Range('D13').Value = "'5.067"
$ws.Range('E13').Value = '  -1.42%  '
$ws.Range('D14').Value = "'0.6830"
$ws.Range('E14').Value = '  -0.41%  '
$ws.Range('D15').Value = "'87.60"
$ws.Range('E15').Value = '  -5.78%  '
$ws.Range('D16').Value = "'6.177"
$ws.Range('E16').Value = '  -6.45%  '
$ws.Range('D17').Value = '29.159.00'
$ws.Range('E17').Value = '  -0.15%  '
$ws.Range('D18').Value = "'0.000008164"
$ws.Range('E18').Value = '  -1.70%  '
$ws.Range('D19').Value = '2.077.36'
$ws.Range('E19').Value = '  -0.92%  '
$ws.Range('D20').Value = "'228.57"
$ws.Range('E20').Value = '  -5.36%  '
$ws.Range('D21').Value = "'12.56"
$ws.Range('E21').Value = '  -1.27%  '
$ws.Range('D22').Value = "'0.9993"
$ws.Range('E22').Value = '  -0.08%  '
$ws.Range('D23').Value = "'7.408"
$ws.Range('E23').Value = '  -1.43%  '
$ws.Range('D24').Value = "'0.9994"
$ws.Range('E24').Value = '  -0.12%  '
$ws.Range('D25').Value = "'0.1456"
$ws.Range('E25').Value = '  -3.65%  '
$ws.Range('D26').Value = "'159.63"
$ws.Range('E26').Value = '  +0.22%  '
$ws.Range('D27').Value = "'8.781"
$ws.Range('E27').Value = '  -0.47%  '
$ws.Range('D28').Value = "'18.10"
$ws.Range('E28').Value = '  -1.03%  '
$ws.Range('D29').Value = "'1.512"
$ws.Range('E29').Value = '  -1.67%  '
$ws.Range('D30').Value = "'4.281"
$ws.Range('E30').Value = '  +1.32%  '
$ws.Range('D31').Value = "'4.141"
$ws.Range('E31').Value = '  -0.91%  '
$ws.Range('E32').Value = '  -0.97%  '
$ws.Range('D33').Value = "'0.05256"
$ws.Range('E33').Value = '  +2.62%  '
$ws.Range('D34').Value = "'0.7582"
$ws.Range('E34').Value = '  -4.14%  '
$ws.Range('D35').Value = "'1.854"
$ws.Range('E35').Value = '  -2.79%  '
$ws.Range('E36').Value = '  -1.07%  '
$ws.Range('E37').Value = '  -0.41%  '
$ws.Range('D38').Value = '1.305.39'
$ws.Range('E38').Value = '  -1.20%  '
$ws.Range('D39').Value = "'0.01832"
$ws.Range('E39').Value = '  -1.76%  '
$ws.Range('D40').Value = "'2.726"
$ws.Range('E40').Value = '  +0.55%  '
$ws.Range('D41').Value = "'0.9298"
$ws.Range('E41').Value = '  -3.28%  '
$ws.Range('D42').Value = "'5.938"
$ws.Range('E42').Value = '  -2.11%  '
$ws.Range('D43').Value = "'104.92"
$ws.Range('E43').Value = '  -2.18%  '
$ws.Range('E44').Value = '  -0.11%  '
$ws.Range('B45').Value = 'Aave'
$ws.Range('C45').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D45').Value = "'65.04"
$ws.Range('E45').Value = '  +1.00%  '
$ws.Range('B46').Value = 'RocketPoolETH'
$ws.Range('C46').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D46').Value = '1.979.96'
$ws.Range('E46').Value = '  -0.64%  '
$ws.Range('B47').Value = 'Mantle'
$ws.Range('C47').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D47').Value = "'0.5196"
$ws.Range('E47').Value = '  +0.24%  '
$ws.Range('E48').Value = '  -0.47%  '
$ws.Range('D49').Value = "'9.548"
$ws.Range('E49').Value = '  -1.59%  '
$ws.Range('D50').Value = "'1.773"
$ws.Range('E50').Value = '  +0.62%  '
$ws.Range('D51').Value = "'0.05951"
$ws.Range('E51').Value = '  +0.87%  '

# Reset style on cells where a text-forcing apostrophe was used,
# so the cell keeps default (unstyled) formatting like the rest of the data cells.
$ws.Range('D5').Style = 'Normal'
$ws.Range('D6').Style = 'Normal'
$ws.Range('D7').Style = 'Normal'
$ws.Range('D8').Style = 'Normal'
$ws.Range('D9').Style = 'Normal'
$ws.Range('D10').Style = 'Normal'
$ws.Range('D11').Style = 'Normal'
$ws.Range('D13').Style = 'Normal'
$ws.Range('D14').Style = 'Normal'
$ws.Range('D15').Style = 'Normal'
$ws.Range('D16').Style = 'Normal'
$ws.Range('D18').Style = 'Normal'
$ws.Range('D20').Style = 'Normal'
$ws.Range('D21').Style = 'Normal'
$ws.Range('D22').Style = 'Normal'
$ws.Range('D23').Style = 'Normal'
$ws.Range('D24').Style = 'Normal'
$ws.Range('D25').Style = 'Normal'
$ws.Range('D26').Style = 'Normal'
$ws.Range('D27').Style = 'Normal'
$ws.Range('D28').Style = 'Normal'
$ws.Range('D29').Style = 'Normal'
$ws.Range('D30').Style = 'Normal'
$ws.Range('D31').Style = 'Normal'
$ws.Range('D33').Style = 'Normal'
$ws.Range('D34').Style = 'Normal'
$ws.Range('D35').Style = 'Normal'
$ws.Range('D39').Style = 'Normal'
$ws.Range('D40').Style = 'Normal'
$ws.Range('D41').Style = 'Normal'
$ws.Range('D42').Style = 'Normal'
$ws.Range('D43').Style = 'Normal'
$ws.Range('D45').Style = 'Normal'
$ws.Range('D47').Style = 'Normal'
$ws.Range('D49').Style = 'Normal'
$ws.Range('D50').Style = 'Normal'
$ws.Range('D51').Style = 'Normal'
